$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Execution Status column (F) for the "Register" scenario rows
# from "No" to "Yes" so UFT One executes the corresponding test case.
$ws.Range("F18").Value = "Yes"
$ws.Range("F19").Value = "Yes"

# Update the selected range/active cell to match the last user selection.
$ws.Range("F24:F25").Select()
